$wb = $excel.ActiveWorkbook

# Sheet 1 ("展览") updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 12420
$ws1.Range("F5").Value = 258
$ws1.Range("F8").Value = 12359
$ws1.Range("F10").Value = 525
$ws1.Range("F17").Value = 6042
$ws1.Range("F20").Value = 212

# Sheet 4 ("全部类型") updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 12420
$ws4.Range("F5").Value = 258
$ws4.Range("F9").Value = 12359
$ws4.Range("F11").Value = 525
$ws4.Range("F19").Value = 6042
$ws4.Range("F22").Value = 212
